$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 data (highlighted in red)
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = "aaa"
$ws.Range("C5").Value = "0.5 mg"
$ws.Range("D5").Value = "Cajita"
$ws.Range("E5").Value = "Amount1"
$ws.Range("F5").Value = "Crema"
$ws.Range("G5").Value = "comerciales"
$ws.Range("I5").Value = 1.2

# Row 6 data
$ws.Range("A6").Value = 2
$ws.Range("B6").Value = "Baselina"
$ws.Range("C6").Value = "600 mg"
$ws.Range("D6").Value = "Caja vial"
$ws.Range("E6").Value = "Novartis"
$ws.Range("F6").Value = "Crema"
$ws.Range("G6").Value = "Sobresitos"
$ws.Range("H6").Value = 132
$ws.Range("I6").Value = 2.7

# Row 7 data
$ws.Range("A7").Value = 3
$ws.Range("B7").Value = "Mentisan"
$ws.Range("C7").Value = "500gr."
$ws.Range("D7").Value = "null"
$ws.Range("E7").Value = "Inti"
$ws.Range("F7").Value = "Crema"
$ws.Range("G7").Value = "comerciales"
$ws.Range("H7").Value = 86
$ws.Range("I7").Value = 7.5

# Apply red font color to row 5 (A5:I5), including the empty H5 cell
$row5 = $ws.Range("A5:I5")
$row5.Font.Color = 255

# Column widths (bestFit/custom) for columns B..I
$ws.Columns.Item(2).ColumnWidth = 9.166666667
$ws.Columns.Item(3).ColumnWidth = 15.166666667
$ws.Columns.Item(4).ColumnWidth = 10.166666667
$ws.Columns.Item(5).ColumnWidth = 12.166666667
$ws.Columns.Item(6).ColumnWidth = 14.166666667
$ws.Columns.Item(7).ColumnWidth = 12.166666667
$ws.Columns.Item(8).ColumnWidth = 5.166666667
$ws.Columns.Item(9).ColumnWidth = 7.166666667

# Selection moves to the newly entered row
$ws.Range("A5:I5").Select() | Out-Null
